# Fix inconsistent "WT " (trailing space) genotype labels so that they
# match the canonical "WT" value already used elsewhere in the workbook.
#
# Affected sheets/columns:
#   - "Pole and Beam" : Genotype values live in column D
#   - "Sheet1"        : Genotype values live in column B
#
# Any cell in those columns whose text is exactly "WT " (with a trailing
# space) is rewritten to "WT".

$wb = $excel.ActiveWorkbook

$targets = @(
    @{ Sheet = "Pole and Beam"; Col = "D" },
    @{ Sheet = "Sheet1";        Col = "B" }
)

foreach ($target in $targets) {
    $ws = $wb.Worksheets.Item($target.Sheet)
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Range($target.Col + $r)
        if ($cell.Value2 -eq "WT ") {
            $cell.Value = "WT"
        }
    }
}
